$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename OTP template codes (rows 5,6,7,11,12,13,17,18,19) ---
$ws.Range("A5").Value = 'ida-auth-otp-email-content-template'
$ws.Range("A6").Value = 'ida-auth-otp-email-subject-template'
$ws.Range("A7").Value = 'ida-auth-otp-sms-template'
$ws.Range("A11").Value = 'ida-auth-otp-email-content-template'
$ws.Range("A12").Value = 'ida-auth-otp-email-subject-template'
$ws.Range("A13").Value = 'ida-auth-otp-sms-template'
$ws.Range("A17").Value = 'ida-auth-otp-email-content-template'
$ws.Range("A18").Value = 'ida-auth-otp-email-subject-template'
$ws.Range("A19").Value = 'ida-auth-otp-sms-template'

# --- Append new rows 125-136 ---
$newRows = New-Object 'object[,]' 12,6
$newRows[0,0] = 'consent'
$newRows[0,1] = 'Consent'
$newRows[0,2] = 'eng'
$newRows[0,3] = $true
$newRows[0,4] = 'superadmin'
$newRows[0,5] = 'now()'
$newRows[1,0] = 'consent'
$newRows[1,1] = 'موافقة'
$newRows[1,2] = 'ara'
$newRows[1,3] = $true
$newRows[1,4] = 'superadmin'
$newRows[1,5] = 'now()'
$newRows[2,0] = 'consent'
$newRows[2,1] = 'Consentement'
$newRows[2,2] = 'fra'
$newRows[2,3] = $true
$newRows[2,4] = 'superadmin'
$newRows[2,5] = 'now()'
$newRows[3,0] = 'auth-otp-email-subject-template'
$newRows[3,1] = 'Auth OTP Email Subject Template'
$newRows[3,2] = 'eng'
$newRows[3,3] = $true
$newRows[3,4] = 'superadmin'
$newRows[3,5] = 'now()'
$newRows[4,0] = 'auth-otp-email-subject-template'
$newRows[4,1] = 'مصادقة OTP قالب موضوع'
$newRows[4,2] = 'ara'
$newRows[4,3] = $true
$newRows[4,4] = 'superadmin'
$newRows[4,5] = 'now()'
$newRows[5,0] = 'auth-otp-email-subject-template'
$newRows[5,1] = 'Modèle dobjet de-mail Auth OTP'
$newRows[5,2] = 'fra'
$newRows[5,3] = $true
$newRows[5,4] = 'superadmin'
$newRows[5,5] = 'now()'
$newRows[6,0] = 'auth-otp-email-content-template'
$newRows[6,1] = 'Auth OTP Email Content Template'
$newRows[6,2] = 'eng'
$newRows[6,3] = $true
$newRows[6,4] = 'superadmin'
$newRows[6,5] = 'now()'
$newRows[7,0] = 'auth-otp-email-content-template'
$newRows[7,1] = 'مصادقة OTP قالب محتوى'
$newRows[7,2] = 'ara'
$newRows[7,3] = $true
$newRows[7,4] = 'superadmin'
$newRows[7,5] = 'now()'
$newRows[8,0] = 'auth-otp-email-content-template'
$newRows[8,1] = 'Auth OTP Email ContentTemplate'
$newRows[8,2] = 'fra'
$newRows[8,3] = $true
$newRows[8,4] = 'superadmin'
$newRows[8,5] = 'now()'
$newRows[9,0] = 'auth-otp-sms-template'
$newRows[9,1] = 'Auth OTP SMS Template'
$newRows[9,2] = 'eng'
$newRows[9,3] = $true
$newRows[9,4] = 'superadmin'
$newRows[9,5] = 'now()'
$newRows[10,0] = 'auth-otp-sms-template'
$newRows[10,1] = 'مصادقة قالب رسالة OTP'
$newRows[10,2] = 'ara'
$newRows[10,3] = $true
$newRows[10,4] = 'superadmin'
$newRows[10,5] = 'now()'
$newRows[11,0] = 'auth-otp-sms-template'
$newRows[11,1] = 'Modèle SMS OTP Auth'
$newRows[11,2] = 'fra'
$newRows[11,3] = $true
$newRows[11,4] = 'superadmin'
$newRows[11,5] = 'now()'

$ws.Range("A125:F136").Value = $newRows

# --- Update selection to reflect new used range (A137 onward) ---
[void]$ws.Range("A137:XFD1048576").Select()
